$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.496.14"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "1.841.78"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "262.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5321"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3089"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06896"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07830"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7603"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.36%  "
$ws.Range("D13").Value = "1.840.34"
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "89.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.045"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007950"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").Value = "26.531.55"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.630"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.317"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "141.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.190"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.690"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "111.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.282"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08810"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.094"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04829"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.934"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7330"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.134"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.103"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.321"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.37%  "
$ws.Range("E38").Value = "  -3.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.4807"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9030"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "108.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.893"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.65%  "
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.634"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.078"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("E47").Value = "  +1.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9005"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05806"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "60.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.89%  "
